$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the amount column values on row 2 (stored as text, same as before the edit)
$ws.Range("C2:D2").NumberFormat = "@"
$ws.Range("C2").Value = "2016-10-15"
$ws.Range("D2").Value = "3324"
$ws.Range("C2:D2").ClearFormats()

# Add a new row with the doubles needed for calculation
$ws.Range("A3").Value = 2.0

$ws.Range("B3:E3").NumberFormat = "@"
$ws.Range("B3").Value = "1"
$ws.Range("C3").Value = "2016-10-15"
$ws.Range("D3").Value = "123"
$ws.Range("E3").Value = "LOL AMOUNT"
$ws.Range("B3:E3").ClearFormats()
